$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "category Foldable2<f> =" -> "category Foldable2<f; Foldable<f>> ="
#    The trailing "<f> =" run is split into three runs:
#      "<f; "  "Foldable<f>"  "> ="
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("category Foldable2<f> =") | Out-Null
$oldEnd = $rng.End
$oldStart = $oldEnd - 5   # "<f> =" is 5 characters long

$piece1 = $d.Range($oldEnd, $oldEnd)
$piece1.InsertAfter("<f; ")
# Stamp this run's Bold so later edits elsewhere in the doc don't coalesce
# it together with its neighbour runs.
$piece1.Font.Bold = $true
$piece1.Font.Bold = $false

$piece2 = $d.Range($piece1.End, $piece1.End)
$piece2.InsertAfter("Foldable<f>")
$piece2.Font.Bold = $true
$piece2.Font.Bold = $false

$piece3 = $d.Range($piece2.End, $piece2.End)
$piece3.InsertAfter("> =")
$piece3.Font.Bold = $true
$piece3.Font.Bold = $false

# Remove the original "<f> =" text that is now left before our new pieces.
$oldRange = $d.Range($oldStart, $oldStart + 5)
$oldRange.Text = ""

# ---------------------------------------------------------------------------
# 2) "| fold2<a, b, c; Foldable<a>; Foldable<b>> : ..." ->
#    "| fold2<a, b, c> : ..."
#    The first two runs merge into a single run; the trailing run (starting
#    at "> : (c -> ...") is left untouched.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("| fold2<a, b, c; Foldable<a>; Foldable<b>") | Out-Null
$rng2.Text = "| fold2<a, b, c"
# Stamp to prevent this run from silently merging with the following,
# identically-formatted "> : (c -> a -> b -> c) -> f<a> -> f<b> -> c" run.
$rng2.Font.Bold = $true
$rng2.Font.Bold = $false

# ---------------------------------------------------------------------------
# 3) Bold the six short "function name" runs that close out the
#    *Stream definitions (map, fold, map2 [as "map" + "2"], product, sum).
# ---------------------------------------------------------------------------
function Set-TrailingBold($precedingText, $boldLen) {
    $r = $d.Content
    $r.Find.Execute($precedingText) | Out-Null
    $target = $d.Range($r.End, $r.End + $boldLen)
    $target.Font.Bold = $true
    $target.Font.BoldBi = $true
}

Set-TrailingBold "mapStream<a, b> (a -> b) -> Stream<a> -> Stream<b> = " 3
Set-TrailingBold "foldStream<a, b> : (b -> a -> b) -> b -> Stream<a> -> b = " 4
Set-TrailingBold "Stream<c>) -> Stream<a> -> Stream<b> -> Stream<c> = " 3
Set-TrailingBold "Stream<a, b, c> (a -> b -> c) -> Stream<a> -> Stream<b> -> Stream<c> = map" 1
Set-TrailingBold "productStream<a, b> : Stream<a> -> Stream<b> -> Stream<(a, b)> = " 7
Set-TrailingBold "sumStream<a, b> : Stream<a> -> Stream<b> -> Stream<Either<a, b>> = " 3

Write-Output "done"
